# FUN_YR_FIN.xlsx update:
#  - A new "most recent period" column is inserted as column D (the old
#    D:K data, i.e. the seven prior periods, shifts right to E:L).
#  - The new column D is populated with the latest period's figures.
#  - A handful of previously-reported prior-period figures are corrected
#    (restated) as part of the same update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new column at D, shifting existing D:K -> E:L ------------
$ws.Range("D:D").Insert()

# The freshly inserted column D has no formatting of its own yet; copy the
# number/font formatting from column E (which used to be D) so the new
# column matches the rest of the table (date format on row 7/38/80, the
# Verdana number format elsewhere, etc).
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match column D's width to its neighbours.
$ws.Range("D1").ColumnWidth = $ws.Range("E1").ColumnWidth

# --- 2. Populate the new column D with the latest-period values ----------
$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(8, 4).Value = 1348500
$ws.Cells.Item(9, 4).Value = 114700
$ws.Cells.Item(10, 4).Value = 1233800
$ws.Cells.Item(12, 4).Value = "NA"
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(14, 4).Value = 11300
$ws.Cells.Item(15, 4).Value = 155500
$ws.Cells.Item(17, 4).Value = 1059100
$ws.Cells.Item(18, 4).Value = 289400
$ws.Cells.Item(20, 4).Value = -42400
$ws.Cells.Item(21, 4).Value = 402600
$ws.Cells.Item(22, 4).Value = 85700
$ws.Cells.Item(23, 4).Value = 161400
$ws.Cells.Item(24, 4).Value = 44600
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(26, 4).Value = 116800
$ws.Cells.Item(27, 4).Value = 116800
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(29, 4).Value = 9900
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(32, 4).Value = 42400
$ws.Cells.Item(33, 4).Value = 126700
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(35, 4).Value = 126700

$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(41, 4).Value = 105300
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(43, 4).Value = 51500
$ws.Cells.Item(44, 4).Value = 30800
$ws.Cells.Item(45, 4).Value = 12600
$ws.Cells.Item(46, 4).Value = 200200
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(48, 4).Value = 1599400
$ws.Cells.Item(49, 4).Value = 215100
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(52, 4).Value = 9400
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(54, 4).Value = 2024200
$ws.Cells.Item(57, 4).Value = 23300
$ws.Cells.Item(58, 4).Value = 5600
$ws.Cells.Item(59, 4).Value = 205800
$ws.Cells.Item(60, 4).Value = 234700
$ws.Cells.Item(61, 4).Value = 1657600
$ws.Cells.Item(62, 4).Value = 99500
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(66, 4).Value = 1991800
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(72, 4).Value = 0
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(76, 4).Value = 32400
$ws.Cells.Item(77, 4).Value = 0

$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(81, 4).Value = 126700
$ws.Cells.Item(83, 4).Value = 155500
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(89, 4).Value = 350700
$ws.Cells.Item(91, 4).Value = -189800
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(94, 4).Value = -189700
$ws.Cells.Item(96, 4).Value = -203200
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(100, 4).Value = -216600
$ws.Cells.Item(101, 4).Value = -5400
$ws.Cells.Item(102, 4).Value = -60900

# --- 3. A few prior-period figures were restated at the same time --------
# (row 89 = "Total Cash Flows From Investing Activities",
#  row 91 = "Capital Expenditures",
#  row 100 = "Change In Cash and Cash Equivalents")
$ws.Cells.Item(89, 6).Value = 358300
$ws.Cells.Item(91, 5).Value = -188100
$ws.Cells.Item(91, 6).Value = -160700
$ws.Cells.Item(100, 6).Value = -194500
